$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- New header cell for the added "vars_internal" column (F) ---
$ws.Range("F1").Value = "vars_internal"

# Copy the header formatting (bold font, border, centered/top alignment)
# from an existing header cell onto the new one instead of rebuilding the
# style by hand, so the shared cellXfs/border/font tables stay untouched.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# --- Simplify/clean the docstring column (E) values ---
$ws.Range("E2").Value = "A docstring for a procedure"
$ws.Range("E3").Value = "Method1 docstring is" + $nl + "multiline"
$ws.Range("E4").Value = "Method2 docstring"

# --- Populate the new vars_internal column (F) with each routine's locals ---
$ws.Range("F2").Value = "i|Integer," + $nl + "j|Integer"
$ws.Range("F3").Value = "cht|Chart|New," + $nl + "tbl|tblRowsCols|New"
$ws.Range("F4").Value = "rng|Variant," + $nl + "tbl|Object"
